# "added proposal script for plan details"
# Adds new proposal-form fields (sumAssured/policyTerm/FreqPayment/PSDay/
# PSMonth/PSYear/MedClass) to the ProposalForm sheet, renumbers the
# proposal number, and re-types the existing sumAssured sample value as
# text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProposalForm")
$ws.Activate()

# --- existing-field edits -------------------------------------------------

# New sample proposal number.
$ws.Range("B2").Value = "LN000002"

# AI2 held the numeric sample 100000 for sumAssured; it becomes a text
# value (leading apostrophe forces text typing, keeping the cell's
# existing left-aligned style).
$ws.Range("AI2").Value = "'100000"

# --- new header row (row 1) ------------------------------------------------

$ws.Range("AS1").Value = "policyTerm"
$ws.Range("AT1").Value = "FreqPayment"
$ws.Range("AU1").Value = "PSDay"
$ws.Range("AV1").Value = "PSMonth"
$ws.Range("AW1").Value = "PSYear"
$ws.Range("AX1").Value = "MedClass"

# Match the bold header formatting used by every other row-1 header cell.
$ws.Range("AS1:AX1").Font.Bold = $true

# --- new sample data row (row 2) -------------------------------------------

# AR2 held the numeric sample 500000; the new sample is 700000, stored as
# text like its neighbouring quote-prefixed cells.
$ws.Range("AR2").Value = "'700000"
$ws.Range("AS2").Value = "'25"
$ws.Range("AT2").Value = "Monthly"
$ws.Range("AU2").Value = "'20"
$ws.Range("AV2").Value = "'4"
$ws.Range("AW2").Value = "'2019"
$ws.Range("AX2").Value = "Non-Medical"

# --- column widths for the newly introduced columns ------------------------

$ws.Columns("AR").ColumnWidth = 10.592447916666666
$ws.Columns("AS").ColumnWidth = 9.877604166666666
$ws.Columns("AT").ColumnWidth = 11.877604166666666
$ws.Columns("AX").ColumnWidth = 11.736979166666666

# --- selection follows the last-edited cell, like the authored edit -------

[void]$ws.Range("AX2").Select()
